# Apply updated crypto price/volume figures from the data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.959.68'
$ws.Range("E2").Value = '  -0.42%  '
$ws.Range("D3").Value = '2.504.81'
$ws.Range("E3").Value = '  -0.70%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").Value = '''534.75'
$ws.Range("D5").ClearFormats()  # drop the quote-prefix text style Excel applied
$ws.Range("E5").Value = '  -0.12%  '
$ws.Range("D6").Value = '''137.11'
$ws.Range("D6").ClearFormats()  # drop the quote-prefix text style Excel applied
$ws.Range("E6").Value = '  -2.77%  '
$ws.Range("D7").Value = '''0.995'
$ws.Range("D7").ClearFormats()  # drop the quote-prefix text style Excel applied
$ws.Range("E7").Value = '  -0.42%  '
$ws.Range("D8").Value = '''0.561'
$ws.Range("D8").ClearFormats()  # drop the quote-prefix text style Excel applied
$ws.Range("E8").Value = '  -0.77%  '
$ws.Range("D9").Value = '2.534.25'
$ws.Range("E9").Value = '  +0.11%  '
$ws.Range("D10").Value = '''0.101'
$ws.Range("D10").ClearFormats()  # drop the quote-prefix text style Excel applied
$ws.Range("E10").Value = '  +1.72%  '
$ws.Range("E11").Value = '  -0.33%  '
$ws.Range("D12").Value = '''5.32'
$ws.Range("D12").ClearFormats()  # drop the quote-prefix text style Excel applied
$ws.Range("E12").Value = '  -2.02%  '
$ws.Range("D13").Value = '''0.348'
$ws.Range("D13").ClearFormats()  # drop the quote-prefix text style Excel applied
$ws.Range("E13").Value = '  -1.88%  '
$ws.Range("D14").Value = '2.946.71'
$ws.Range("E14").Value = '  -1.00%  '
$ws.Range("D15").Value = '''23.19'
$ws.Range("D15").ClearFormats()  # drop the quote-prefix text style Excel applied
$ws.Range("E15").Value = '  -0.60%  '
$ws.Range("D16").Value = '58.840.95'
$ws.Range("E16").Value = '  -0.59%  '
$ws.Range("E17").Value = '  -0.87%  '
$ws.Range("D18").Value = '2.511.80'
$ws.Range("E18").Value = '  -0.40%  '
$ws.Range("D19").Value = '''11.05'
$ws.Range("D19").ClearFormats()  # drop the quote-prefix text style Excel applied
$ws.Range("E19").Value = '  +0.25%  '
$ws.Range("D20").Value = '''4.25'
$ws.Range("D20").ClearFormats()  # drop the quote-prefix text style Excel applied
$ws.Range("E20").Value = '  +0.18%  '
$ws.Range("D21").Value = '''324.59'
$ws.Range("D21").ClearFormats()  # drop the quote-prefix text style Excel applied
$ws.Range("E21").Value = '  +0.97%  '
$ws.Range("D22").Value = '''0.999'
$ws.Range("D22").ClearFormats()  # drop the quote-prefix text style Excel applied
$ws.Range("E22").Value = '  +0.08%  '
$ws.Range("D23").Value = '''5.88'
$ws.Range("D23").ClearFormats()  # drop the quote-prefix text style Excel applied
$ws.Range("E23").Value = '  +0.92%  '
$ws.Range("D24").Value = '''63.09'
$ws.Range("D24").ClearFormats()  # drop the quote-prefix text style Excel applied
$ws.Range("E24").Value = '  +0.59%  '
$ws.Range("D25").Value = '''0.420'
$ws.Range("D25").ClearFormats()  # drop the quote-prefix text style Excel applied
$ws.Range("E25").Value = '  -0.50%  '
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("D27").Value = '''0.992'
$ws.Range("D27").ClearFormats()  # drop the quote-prefix text style Excel applied
$ws.Range("E27").Value = '  -1.08%  '
$ws.Range("D28").Value = '''7.59'
$ws.Range("D28").ClearFormats()  # drop the quote-prefix text style Excel applied
$ws.Range("E28").Value = '  -3.24%  '
$ws.Range("D29").Value = '''6.85'
$ws.Range("D29").ClearFormats()  # drop the quote-prefix text style Excel applied
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("D30").Value = '''0.0₃0777'
$ws.Range("D30").ClearFormats()  # drop the quote-prefix text style Excel applied
$ws.Range("E30").Value = '  +0.48%  '
$ws.Range("E31").Value = '  -0.99%  '
$ws.Range("D32").Value = '''166.74'
$ws.Range("D32").ClearFormats()  # drop the quote-prefix text style Excel applied
$ws.Range("E32").Value = '  +1.57%  '
$ws.Range("B33").Value = 'Fetch.AI'
$ws.Range("C33").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D33").Value = '''1.14'
$ws.Range("D33").ClearFormats()  # drop the quote-prefix text style Excel applied
$ws.Range("E33").Value = '  +0.23%  '
$ws.Range("B34").Value = 'USDe'
$ws.Range("C34").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D34").Value = '''0.997'
$ws.Range("D34").ClearFormats()  # drop the quote-prefix text style Excel applied
$ws.Range("E34").Value = '  -0.17%  '
$ws.Range("D35").Value = '''1.40'
$ws.Range("D35").ClearFormats()  # drop the quote-prefix text style Excel applied
$ws.Range("E35").Value = '  -2.29%  '
$ws.Range("D36").Value = '''18.49'
$ws.Range("D36").ClearFormats()  # drop the quote-prefix text style Excel applied
$ws.Range("E36").Value = '  -0.08%  '
$ws.Range("D37").Value = '''4.13'
$ws.Range("D37").ClearFormats()  # drop the quote-prefix text style Excel applied
$ws.Range("E37").Value = '  -3.75%  '
$ws.Range("E38").Value = '  -1.93%  '
$ws.Range("D39").Value = '''36.72'
$ws.Range("D39").ClearFormats()  # drop the quote-prefix text style Excel applied
$ws.Range("E39").Value = '  -0.41%  '
$ws.Range("D40").Value = '''0.833'
$ws.Range("D40").ClearFormats()  # drop the quote-prefix text style Excel applied
$ws.Range("E40").Value = '  +2.97%  '
$ws.Range("D41").Value = '''3.62'
$ws.Range("D41").ClearFormats()  # drop the quote-prefix text style Excel applied
$ws.Range("E41").Value = '  -0.93%  '
$ws.Range("D42").Value = '''5.27'
$ws.Range("D42").ClearFormats()  # drop the quote-prefix text style Excel applied
$ws.Range("E42").Value = '  -2.04%  '
$ws.Range("D43").Value = '''279.53'
$ws.Range("D43").ClearFormats()  # drop the quote-prefix text style Excel applied
$ws.Range("E43").Value = '  -3.64%  '
$ws.Range("E44").Value = '  -0.31%  '
$ws.Range("D45").Value = '''0.606'
$ws.Range("D45").ClearFormats()  # drop the quote-prefix text style Excel applied
$ws.Range("E45").Value = '  +1.14%  '
$ws.Range("E46").Value = '  -0.09%  '
$ws.Range("D47").Value = '''125.38'
$ws.Range("D47").ClearFormats()  # drop the quote-prefix text style Excel applied
$ws.Range("E47").Value = '  +0.32%  '
$ws.Range("D48").Value = '''0.0928'
$ws.Range("D48").ClearFormats()  # drop the quote-prefix text style Excel applied
$ws.Range("E48").Value = '  -0.24%  '
$ws.Range("D49").Value = '''0.0511'
$ws.Range("D49").ClearFormats()  # drop the quote-prefix text style Excel applied
$ws.Range("E49").Value = '  +0.02%  '
$ws.Range("E50").Value = '  -1.20%  '
$ws.Range("D51").Value = '''17.58'
$ws.Range("D51").ClearFormats()  # drop the quote-prefix text style Excel applied
$ws.Range("E51").Value = '  -0.33%  '
